$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rows 4 and 5: add an empty H cell styled like the existing F66-style cells (cellXfs index 1) ---
$ws.Range("F66").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("H5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Rows 66-85: Temperature_C (G) and pH_total_scale (H) readings ---
$data = @{
    66 = @(15,    7.86)
    67 = @(15.2,  7.6)
    68 = @(9.2,   7.78)
    69 = @(8.9,   7.65)
    70 = @(9.8,   7.55)
    71 = @(15.1,  7.59)
    72 = @(15.2,  7.74)
    73 = @(9.8,   7.65)
    74 = @(15.1,  7.58)
    75 = @(9.7,   7.69)
    76 = @(9.1,   7.61)
    77 = @(15,    7.77)
    78 = @(15.1,  7.75)
    79 = @(15,    7.7)
    80 = @(9.4,   7.71)
    81 = @(15.3,  7.62)
    82 = @(14.9,  7.6)
    83 = @(9.3,   7.6)
    84 = @(9.3,   7.68)
    85 = @(9.4,   7.66)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $gCell = $ws.Cells.Item($row, 7)
    $hCell = $ws.Cells.Item($row, 8)
    $gCell.Value = $vals[0]
    $gCell.Style = "Normal"
    $hCell.Value = $vals[1]
    $hCell.Style = "Normal"
}

# --- Update sheet view (scroll position / selection) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 70
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G77").Select()
